$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name (B) and Link (C) updates (row reordering) ---
$ws.Range("B34").Value = "OKB"
$ws.Range("B35").Value = "TheGraph"
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("B48").Value = "Stacks"
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C35").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

# --- Price (D) updates: force text format so numeric-looking strings
#     (e.g. "1.20", "0.130") keep their exact text representation, ---
#     then clear the format so the cell style index is unaffected. ---
$dCells = @("D2","D3","D4","D5","D6","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D2").Value = "69.413.64"
$ws.Range("D3").Value = "3.948.10"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "502.68"
$ws.Range("D6").Value = "148.34"
$ws.Range("D9").Value = "0.738"
$ws.Range("D10").Value = "0.177"
$ws.Range("D11").Value = "0.0000352"
$ws.Range("D12").Value = "43.86"
$ws.Range("D13").Value = "10.55"
$ws.Range("D14").Value = "4.577.98"
$ws.Range("D15").Value = "3.940.92"
$ws.Range("D16").Value = "14.31"
$ws.Range("D18").Value = "1.20"
$ws.Range("D19").Value = "20.06"
$ws.Range("D20").Value = "69.412.91"
$ws.Range("D21").Value = "437.84"
$ws.Range("D22").Value = "3.45"
$ws.Range("D23").Value = "14.73"
$ws.Range("D24").Value = "89.07"
$ws.Range("D25").Value = "12.03"
$ws.Range("D26").Value = "3.88"
$ws.Range("D28").Value = "37.21"
$ws.Range("D29").Value = "5.67"
$ws.Range("D30").Value = "704.97"
$ws.Range("D31").Value = "13.47"
$ws.Range("D32").Value = "0.130"
$ws.Range("D33").Value = "2.90"
$ws.Range("D34").Value = "64.77"
$ws.Range("D35").Value = "0.456"
$ws.Range("D36").Value = "0.0₃0903"
$ws.Range("D37").Value = "41.30"
$ws.Range("D38").Value = "6.05"
$ws.Range("D39").Value = "0.152"
$ws.Range("D40").Value = "0.997"
$ws.Range("D42").Value = "0.0493"
$ws.Range("D43").Value = "2.89"
$ws.Range("D44").Value = "3.08"
$ws.Range("D45").Value = "3.05"
$ws.Range("D48").Value = "3.01"
$ws.Range("D49").Value = "3.41"
$ws.Range("D50").Value = "0.0₆0346"
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}

# --- Volume(1h) (E) updates ---
$eCells = @("E2","E3","E4","E5","E6","E7","E8","E9","E10","E11","E12","E13","E14","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25","E26","E27","E28","E29","E30","E31","E32","E33","E34","E35","E36","E37","E38","E39","E40","E42","E43","E44","E45","E46","E47","E48","E49","E50","E51")
foreach ($addr in $eCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +5.59%  "
$ws.Range("E26").Value = "  +6.99%  "
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("E35").Value = "  +14.88%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  +7.60%  "
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("E51").Value = "  -2.16%  "
foreach ($addr in $eCells) {
    $ws.Range($addr).ClearFormats()
}
